$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.480.22'
$ws.Range('E2').Value = '  -4.30%  '
$ws.Range('D3').Value = '3.082.57'
$ws.Range('E3').Value = '  -5.53%  '
$ws.Range('E4').Value = '  +0.18%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '607.50'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.02%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '144.40'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -8.46%  '
$ws.Range('E7').Value = '  +0.22%  '
$ws.Range('D8').Value = '3.080.53'
$ws.Range('E8').Value = '  -5.60%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.518'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -4.98%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.148'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -8.27%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '5.19'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -10.41%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.464'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -6.40%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.0000246'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  -9.16%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '34.89'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -10.83%  '
$ws.Range('D15').Value = '3.618.39'
$ws.Range('E15').Value = '  -4.68%  '
$ws.Range('E16').Value = '  +0.99%  '
$ws.Range('D17').Value = '63.607.81'
$ws.Range('E17').Value = '  -4.25%  '
$ws.Range('D18').Value = '3.100.94'
$ws.Range('E18').Value = '  -5.01%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '6.77'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -9.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '471.37'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -6.73%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '14.50'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -6.28%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.694'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -8.18%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.63'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -5.86%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '13.41'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -8.60%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '83.18'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -4.82%  '
$ws.Range('E26').Value = '  -0.11%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '2.77'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  -8.96%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '8.30'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -9.47%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '2.13'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -10.84%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '6.68'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -5.92%  '
$ws.Range('B31').Value = 'FirstDigitalUSD'
$ws.Range('C31').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '1.00'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +0.03%  '
$ws.Range('B32').Value = 'Stacks'
$ws.Range('C32').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.71'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.42%  '
$ws.Range('E33').Value = '  -16.01%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '25.92'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -7.48%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '1.09'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -5.27%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '5.87'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -9.43%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '52.35'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -6.20%  '
$ws.Range('B38').Value = 'Bittensor'
$ws.Range('C38').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '455.84'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -8.18%  '
$ws.Range('B39').Value = 'PEPE'
$ws.Range('C39').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D39').Value = '0.0₃0717'
$ws.Range('E39').Value = '  -8.73%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.87'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -15.15%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.0389'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -7.91%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.118'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -8.38%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '8.27'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -6.56%  '
$ws.Range('D44').Value = '2.809.13'
$ws.Range('E44').Value = '  -6.69%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.262'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -10.68%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '2.22'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -11.98%  '
$ws.Range('E47').Value = '  -0.02%  '
$ws.Range('E48').Value = '  -5.70%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '25.87'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -10.74%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.112'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -5.97%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '117.61'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -2.66%  '

Write-Host "Applied cryptos update"
